$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1:C286").AutoFilter(3, @("0"), 7)
Write-Host "A2 hidden:" $ws.Rows.Item(2).Hidden
Write-Host "A80 hidden:" $ws.Rows.Item(80).Hidden
